$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'297.40"
$ws.Range("E2").Value = "'2.68%"
$ws.Range("E3").Value = "'2.77%"
$ws.Range("D4").Value = "'5.036"
$ws.Range("E4").Value = "'-0.33%"
$ws.Range("D5").Value = "'0.07554"
$ws.Range("E5").Value = "'3.66%"
$ws.Range("D6").Value = "'1.592"
$ws.Range("E6").Value = "'1.94%"
$ws.Range("D7").Value = "'0.9269"
$ws.Range("E7").Value = "'0.70%"
$ws.Range("E9").Value = "'4.66%"
$ws.Range("D10").Value = "'0.1834"
$ws.Range("E10").Value = "'6.47%"
$ws.Range("D11").Value = "'0.08871"
$ws.Range("E11").Value = "'2.29%"
$ws.Range("D12").Value = "'0.03914"
$ws.Range("E12").Value = "'-6.35%"
$ws.Range("D13").Value = "'0.1052"
$ws.Range("E13").Value = "'-0.14%"
$ws.Range("D14").Value = "'0.001278"
$ws.Range("E14").Value = "'0.43%"
$ws.Range("D15").Value = "'0.005770"
$ws.Range("E15").Value = "'-2.24%"
$ws.Range("D16").Value = "'3.333"
$ws.Range("E16").Value = "'-1.87%"
$ws.Range("D17").Value = "'4.375"
$ws.Range("E17").Value = "'2.23%"
$ws.Range("D18").Value = "'0.3321"
$ws.Range("E18").Value = "'1.29%"
$ws.Range("D19").Value = "'7.932"
$ws.Range("E19").Value = "'0.62%"
$ws.Range("D20").Value = "'0.1421"
$ws.Range("E20").Value = "'5.18%"
$ws.Range("E21").Value = "'4.01%"
$ws.Range("D22").Value = "'0.04060"
$ws.Range("E22").Value = "'5.09%"
$ws.Range("D23").Value = "'0.001265"
$ws.Range("E23").Value = "'-0.40%"
$ws.Range("D24").Value = "'0.003994"
$ws.Range("E24").Value = "'3.87%"
$ws.Range("D25").Value = "'0.0001229"
$ws.Range("E25").Value = "'-4.13%"
$ws.Range("E26").Value = "'-0.10%"
$ws.Range("D38").Value = "'0.02414"
$ws.Range("E38").Value = "'4.20%"
$ws.Range("D39").Value = "'0.05207"
$ws.Range("E39").Value = "'4.70%"
$ws.Range("D40").Value = "'0.006385"
$ws.Range("E40").Value = "'-3.85%"
$ws.Range("D41").Value = "'0.007781"
$ws.Range("E41").Value = "'1.32%"
$ws.Range("D42").Value = "'0.1328"
$ws.Range("E42").Value = "'4.34%"
$ws.Range("D43").Value = "'0.007579"
$ws.Range("E43").Value = "'2.80%"
$ws.Range("D44").Value = "'0.007828"
$ws.Range("E44").Value = "'10.78%"
$ws.Range("D45").Value = "'0.3224"
$ws.Range("E45").Value = "'11.45%"
$ws.Range("D46").Value = "'0.00006776"
$ws.Range("E46").Value = "'5.62%"
$ws.Range("D47").Value = "'0.00000000749"
$ws.Range("E47").Value = "'-0.23%"
$ws.Range("D48").Value = "'0.04593"
$ws.Range("E48").Value = "'135.17%"
$ws.Range("D49").Value = "'0.004197"
$ws.Range("E49").Value = "'-0.12%"
$ws.Range("D50").Value = "'0.00002098"
$ws.Range("E50").Value = "'-0.23%"
$ws.Range("D51").Value = "'0.0001998"
$ws.Range("E51").Value = "'-0.23%"
